# Updates cryptos list values (price / 1h volume change) and reorders a few
# coin rows, matching the upstream GitHub Actions scrape commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "20.564.02"
$ws.Range("E2").Value = "  +1.80%  "

# Row 3
$ws.Range("D3").Value = "1.473.36"
$ws.Range("E3").Value = "  +2.95%  "

# Row 4
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.49%  "

# Row 5
$ws.Range("D5").Value = "'0.9673"
$ws.Range("E5").Value = "  -2.81%  "

# Row 6
$ws.Range("D6").Value = "'276.79"
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").Value = "'0.3654"
$ws.Range("E7").Value = "  -1.57%  "

# Row 8
$ws.Range("D8").Value = "'0.3060"

# Row 9
$ws.Range("D9").Value = "'40.65"
$ws.Range("E9").Value = "  +0.70%  "

# Row 10
$ws.Range("D10").Value = "'1.062"
$ws.Range("E10").Value = "  +0.19%  "

# Row 11
$ws.Range("D11").Value = "'0.06641"
$ws.Range("E11").Value = "  +0.72%  "

# Row 12
$ws.Range("D12").Value = "'0.9983"
$ws.Range("E12").Value = "  +0.10%  "

# Row 13
$ws.Range("D13").Value = "'5.477"
$ws.Range("E13").Value = "  -1.59%  "

# Row 14
$ws.Range("D14").Value = "'18.16"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15
$ws.Range("D15").Value = "'6.183"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16
$ws.Range("D16").Value = "'0.00001031"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").Value = "1.475.06"
$ws.Range("E17").Value = "  +2.96%  "

# Row 18
$ws.Range("D18").Value = "'0.9646"
$ws.Range("E18").Value = "  -3.00%  "

# Row 19
$ws.Range("D19").Value = "'0.05912"
$ws.Range("E19").Value = "  +2.72%  "

# Row 20
$ws.Range("D20").Value = "'69.42"
$ws.Range("E20").Value = "  -3.48%  "

# Row 21
$ws.Range("D21").Value = "'5.463"
$ws.Range("E21").Value = "  -3.02%  "

# Row 22
$ws.Range("D22").Value = "'14.55"
$ws.Range("E22").Value = "  -2.22%  "

# Row 23
$ws.Range("D23").Value = "'11.07"
$ws.Range("E23").Value = "  -0.80%  "

# Row 24
$ws.Range("D24").Value = "'2.250"
$ws.Range("E24").Value = "  +0.98%  "

# Row 25
$ws.Range("D25").Value = "20.593.50"
$ws.Range("E25").Value = "  +1.87%  "

# Row 26
$ws.Range("D26").Value = "'140.86"
$ws.Range("E26").Value = "  +4.08%  "

# Row 27
$ws.Range("D27").Value = "'2.140"
$ws.Range("E27").Value = "  -7.57%  "

# Row 28
$ws.Range("D28").Value = "'17.28"
$ws.Range("E28").Value = "  -1.10%  "

# Row 29
$ws.Range("D29").Value = "1.631.51"
$ws.Range("E29").Value = "  +2.40%  "

# Row 30
$ws.Range("D30").Value = "'113.88"
$ws.Range("E30").Value = "  +1.92%  "

# Row 31
$ws.Range("D31").Value = "'3.934"
$ws.Range("E31").Value = "  -0.66%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "'0.8164"
$ws.Range("E32").Value = "  -3.47%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.972"
$ws.Range("E33").Value = "  -6.73%  "

# Row 34
$ws.Range("D34").Value = "'0.07936"
$ws.Range("E34").Value = "  +1.62%  "

# Row 35
$ws.Range("E35").Value = "  +2.72%  "

# Row 36
$ws.Range("D36").Value = "'1.217"
$ws.Range("E36").Value = "  +9.01%  "

# Row 37
$ws.Range("D37").Value = "'0.05825"
$ws.Range("E37").Value = "  -1.77%  "

# Row 38
$ws.Range("D38").Value = "'4.732"
$ws.Range("E38").Value = "  -4.26%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02040"
$ws.Range("E39").Value = "  -1.44%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "'10.48"
$ws.Range("E40").Value = "  -3.21%  "

# Row 41
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.659"
$ws.Range("E41").Value = "  -2.12%  "

# Row 42
$ws.Range("D42").Value = "'0.9648"
$ws.Range("E42").Value = "  -3.05%  "

# Row 43
$ws.Range("D43").Value = "'0.1883"
$ws.Range("E43").Value = "  -0.13%  "

# Row 44
$ws.Range("D44").Value = "'0.5302"
$ws.Range("E44").Value = "  -1.56%  "

# Row 45
$ws.Range("E45").Value = "  -1.60%  "

# Row 46
$ws.Range("D46").Value = "'12.12"
$ws.Range("E46").Value = "  -2.87%  "

# Row 47
$ws.Range("D47").Value = "'118.19"
$ws.Range("E47").Value = "  -0.98%  "

# Row 48
$ws.Range("D48").Value = "'0.5205"
$ws.Range("E48").Value = "  -1.64%  "

# Row 49
$ws.Range("D49").Value = "'1.794"
$ws.Range("E49").Value = "  -0.75%  "

# Row 50
$ws.Range("E50").Value = "  +2.70%  "

# Row 51
$ws.Range("D51").Value = "'0.9926"
$ws.Range("E51").Value = "  -0.50%  "
